$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D3").Value = "40 / 48"
$ws.Range("D6").Value = "55 / 62"
$ws.Range("D7").Value = "88 / 61"

$ws.Range("D8").Select() | Out-Null
